$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.143.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.64%  "
$ws.Range("D3").Value = "'3.052.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'550.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.07%  "
$ws.Range("D6").Value = "'138.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.11%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'3.046.51"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("E9").Value = "  +1.93%  "
$ws.Range("D10").Value = "'6.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.37%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("E13").Value = "  +3.09%  "
$ws.Range("D14").Value = "'34.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.22%  "
$ws.Range("D15").Value = "'3.553.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.60%  "
$ws.Range("D16").Value = "'63.261.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.95%  "
$ws.Range("D17").Value = "'3.054.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.29%  "
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("E19").Value = "  +2.78%  "
$ws.Range("D20").Value = "'479.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.88%  "
$ws.Range("D21").Value = "'13.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.44%  "
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").Value = "'7.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.37%  "
$ws.Range("D24").Value = "'80.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.69%  "
$ws.Range("D25").Value = "'12.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.83%  "
$ws.Range("E27").Value = "  +3.70%  "
$ws.Range("D28").Value = "'7.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.40%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("E30").Value = "  +6.25%  "
$ws.Range("D31").Value = "'25.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.27%  "
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("E33").Value = "  +6.75%  "
$ws.Range("D34").Value = "'5.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.29%  "
$ws.Range("D35").Value = "'55.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("D36").Value = "'5.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.63%  "
$ws.Range("D37").Value = "'461.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.24%  "
$ws.Range("D38").Value = "'0.0809"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.84%  "
$ws.Range("D39").Value = "'3.109.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("E40").Value = "  +2.89%  "
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("D42").Value = "'8.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.74%  "
$ws.Range("E43").Value = "  +4.37%  "
$ws.Range("D44").Value = "'27.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.86%  "
$ws.Range("E45").Value = "  +2.70%  "
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("E47").Value = "  +3.45%  "
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("D49").Value = "'116.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("D50").Value = "'0.0₃0504"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.36%  "
$ws.Range("E51").Value = "  +4.19%  "
